# "Generate Report for handback"
#
# For both locale sheets (zh-cn, de-de) this:
#   - flips the Status column (B) from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - fills in the "Latest Target File" (E) and "Latest Handback File" (F)
#     columns with hyperlinked file names (previously blank)
#   - stamps the "Latest Handback DateTime" (G) with a real timestamp
#     (previously the zero-date placeholder)
#   - flips "Handoff Reason" (H) from "Ignored" to "Include"
#
# Row 4 (.localization-config) is untouched - it is not localizable content
# so it never gets a handback.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# The "Overview" sheet's B2/C2/B3/C3 cells share the same underlying string
# as the per-locale sheets' Status column ("Ready for handoff"). Updating
# that text needs to land here too so every cell that used to read
# "Ready for handoff" now reads the new status.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

function Update-LocaleSheet {
    param(
        [string]$sheetName,
        [string]$xlfFileName,
        [string]$handbackDateTime,
        [string]$commitSha,
        [string]$handoffCommitSha
    )

    $ws = $wb.Worksheets.Item($sheetName)
    Write-Output "Updating sheet: $($ws.Name)"

    $mdFile = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.md"

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSha/e2e/$mdFile"
    $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$handoffCommitSha/ol-handback/OpenLocalizationTestOrg/oltest.$sheetName/xinjiang/$xlfFileName"

    # ---- Row 2 (6ef09566-....md) ----
    $ws.Range("B2").Value = $statusHandedBack

    $ws.Range("E2").Value = $mdFile
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl, "", "", $mdFile) | Out-Null
    $ws.Range("E2").Font.Color = 15570276
    $ws.Range("E2").Font.Underline = 2

    $ws.Range("F2").Value = $xlfFileName
    $ws.Hyperlinks.Add($ws.Range("F2"), $xlfUrl, "", "", $xlfFileName) | Out-Null
    $ws.Range("F2").Font.Color = 15570276
    $ws.Range("F2").Font.Underline = 2

    $ws.Range("G2").Value = $handbackDateTime
    $ws.Range("H2").Value = "Include"

    # ---- Row 3 (ffff868d2531-....md) ----
    $ws.Range("B3").Value = $statusHandedBack

    $ws.Range("E3").Value = $mdFile
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl, "", "", $mdFile) | Out-Null
    $ws.Range("E3").Font.Color = 15570276
    $ws.Range("E3").Font.Underline = 2

    $ws.Range("F3").Value = $xlfFileName
    $ws.Hyperlinks.Add($ws.Range("F3"), $xlfUrl, "", "", $xlfFileName) | Out-Null
    $ws.Range("F3").Font.Color = 15570276
    $ws.Range("F3").Font.Underline = 2

    $ws.Range("G3").Value = $handbackDateTime
    $ws.Range("H3").Value = "Include"
}

Update-LocaleSheet `
    "zh-cn" `
    "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.zh-cn.xlf" `
    "2016-01-19 07:16:06" `
    "6c599477ff7c7d139e3dee8781316052e0909d19" `
    "2efb08e9675d2bf9d7b7a5f73b1fba74de90d02d"

Update-LocaleSheet `
    "de-de" `
    "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.de-de.xlf" `
    "2016-01-19 07:16:23" `
    "6c599477ff7c7d139e3dee8781316052e0909d19" `
    "769e26fbd0e64e6242094d7815ce5ae8ad9bcc16"
